# Atualizado por script em 12-11-2023 20:45
#
# 1) Four pairs of adjacent match rows had their match-data columns (F:V)
#    swapped (the index/country/tournament/season/date columns A:E stayed
#    put) -- this reorders the matches within their already-sorted block
#    while keeping the running index (column A) sequential.
# 2) Three brand-new match rows were appended at the bottom (97-99).
#
# NOTE: this engine's PowerShell function calls only bind parameters
# positionally -- named arguments (`-Foo bar`) are silently dropped, so
# every helper below takes plain positional params.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param($RowA, $RowB)
    $rangeA = $ws.Range("F$RowA`:V$RowA")
    $rangeB = $ws.Range("F$RowB`:V$RowB")
    $tmp = $rangeA.Value2
    $rangeA.Value = $rangeB.Value2
    $rangeB.Value = $tmp
}

Swap-MatchRows 70 71
Swap-MatchRows 76 77
Swap-MatchRows 87 88
Swap-MatchRows 94 95

function Add-MatchRow {
    param($Row, $Indice, $DataPartida, $Home, $HomeGols, $Away, $AwayGols,
          $HomeOpenOdds, $HomeOpenDt, $HomeCloseOdds, $HomeCloseDt,
          $DrawOpenOdds, $DrawOpenDt, $DrawCloseOdds, $DrawCloseDt,
          $AwayOpenOdds, $AwayOpenDt, $AwayCloseOdds, $AwayCloseDt, $Url)

    # Clone formatting (number formats / styles) from the previous row, then
    # overwrite the values so the new row matches the existing table style.
    $prevRow = $Row - 1
    $ws.Range("A$prevRow`:V$prevRow").Copy($ws.Range("A$Row"))

    $ws.Cells.Item($Row, 1).Value = $Indice
    $ws.Cells.Item($Row, 2).Value = "portugal"
    $ws.Cells.Item($Row, 3).Value = "liga-portugal"
    $ws.Cells.Item($Row, 4).Value = "2023-2024"
    $ws.Cells.Item($Row, 5).Value = $DataPartida
    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGols
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGols
    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenDt
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseDt
    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenDt
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseDt
    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenDt
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseDt
    $ws.Cells.Item($Row, 22).Value = $Url
}

Add-MatchRow 97 96 45242.6875 "Boavista" 1 "SC Farense" 3 `
    2.13 "06/11/2023 21:42" 1.97 "12/11/2023 16:14" `
    3.6 "06/11/2023 21:42" 3.73 "12/11/2023 16:14" `
    3.43 "06/11/2023 21:42" 3.93 "12/11/2023 16:14" `
    "https://www.betexplorer.com/football/portugal/liga-portugal/boavista-sc-farense/tCRCo3cH/"

Add-MatchRow 98 97 45242.6875 "Gil Vicente" 1 "Rio Ave" 1 `
    2.11 "05/11/2023 19:12" 2.3 "12/11/2023 16:22" `
    3.56 "05/11/2023 19:12" 3.45 "12/11/2023 16:26" `
    3.61 "05/11/2023 19:12" 3.29 "12/11/2023 16:22" `
    "https://www.betexplorer.com/football/portugal/liga-portugal/gil-vicente-rio-ave/K2QGpNCN/"

Add-MatchRow 99 98 45242.79166666666 "Arouca" 0 "Braga" 1 `
    4.92 "06/11/2023 21:42" 4.75 "12/11/2023 18:58" `
    4.23 "06/11/2023 21:42" 4.02 "12/11/2023 18:58" `
    1.68 "06/11/2023 21:42" 1.75 "12/11/2023 18:58" `
    "https://www.betexplorer.com/football/portugal/liga-portugal/arouca-braga/ppY3mPS4/"
